$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1: the old "Mission" column (J) becomes "Revenue", and a new
# "Mission Statement" column is added in K.
$ws1.Range("J1").Value = "Revenue"
$ws1.Range("K1").Value = "Mission Statement"

# Column J used to be sized for long mission text; now it holds a
# small revenue figure, so narrow it back down.
$ws1.Columns.Item(10).ColumnWidth = 6.83

# Add the new "Skipped" sheet right after Sheet1.
$wsSkipped = $wb.Worksheets.Add($null, $ws1)
$wsSkipped.Name = "Skipped"

$wsSkipped.Range("A1").Value = "Skipped EIN"
$wsSkipped.Range("B1").Value = "Name from IRS Spreadsheet"
$wsSkipped.Range("C1").Value = "Guidestar Link"

# Restore the new sheet's page margins to Excel's out-of-the-box
# defaults (a fresh sheet otherwise inherits the workbook's existing
# non-default margins).
$wsSkipped.PageSetup.LeftMargin = 54
$wsSkipped.PageSetup.RightMargin = 54
$wsSkipped.PageSetup.TopMargin = 72
$wsSkipped.PageSetup.BottomMargin = 72
$wsSkipped.PageSetup.HeaderMargin = 36
$wsSkipped.PageSetup.FooterMargin = 36

$ws1.Activate()
